# Auto-generated edit script: updates numeric cell values per the commit diff
# (profit/price recalculation refresh across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 177.57
$ws.Range("I15").Value = 177.57
$ws.Range("K15").Value = 532.71
$ws.Range("M15").Value = -363.71
$ws.Range("H132").Value = 3091.2058
$ws.Range("I132").Value = 3338.2
$ws.Range("J132").Value = 1238.75
$ws.Range("K132").Value = 10014.6
$ws.Range("L132").Value = 3716.25
$ws.Range("M132").Value = -7484.599999999999
$ws.Range("N132").Value = -8776.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10015.775
$ws.Range("I32").Value = 6425.2705
$ws.Range("K32").Value = 6425.2705
$ws.Range("M32").Value = -6138.2705
$ws.Range("H61").Value = 3370065.5
$ws.Range("I61").Value = 4118506.8
$ws.Range("J61").Value = 2080
$ws.Range("K61").Value = 4118506.8
$ws.Range("L61").Value = 2080
$ws.Range("M61").Value = -4118294.8
$ws.Range("N61").Value = -2504
$ws.Range("H74").Value = 23817276
$ws.Range("I74").Value = 38462556
$ws.Range("J74").Value = 18694.625
$ws.Range("K74").Value = 38462556
$ws.Range("L74").Value = 18694.625
$ws.Range("M74").Value = -38461682
$ws.Range("N74").Value = -20442.625
$ws.Range("H77").Value = 23817276
$ws.Range("I77").Value = 38462556
$ws.Range("J77").Value = 18694.625
$ws.Range("K77").Value = 192312780
$ws.Range("L77").Value = 93473.125
$ws.Range("M77").Value = -192308412
$ws.Range("N77").Value = -102209.125
$ws.Range("H102").Value = 1765.9333
$ws.Range("I102").Value = 1459.0834
$ws.Range("J102").Value = 2993.3333
$ws.Range("K102").Value = 1459.0834
$ws.Range("L102").Value = 2993.3333
$ws.Range("M102").Value = 162.9166
$ws.Range("N102").Value = -6237.3333
$ws.Range("H122").Value = 2091.625
$ws.Range("I122").Value = 1944
$ws.Range("J122").Value = 2239.25
$ws.Range("K122").Value = 5832
$ws.Range("L122").Value = 6717.75
$ws.Range("M122").Value = -3382
$ws.Range("N122").Value = -11617.75
$ws.Range("H132").Value = 988729.5
$ws.Range("I132").Value = 1234582.2
$ws.Range("J132").Value = 128245
$ws.Range("K132").Value = 3703746.6
$ws.Range("L132").Value = 384735
$ws.Range("M132").Value = -3701216.6
$ws.Range("N132").Value = -389795
$ws.Range("H136").Value = 3370065.5
$ws.Range("I136").Value = 4118506.8
$ws.Range("J136").Value = 2080
$ws.Range("K136").Value = 12355520.4
$ws.Range("L136").Value = 6240
$ws.Range("M136").Value = -12352970.4
$ws.Range("N136").Value = -11340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 70400.60000000001
$ws.Range("J57").Value = 70400.60000000001
$ws.Range("L57").Value = 70400.60000000001
$ws.Range("N57").Value = -71840.60000000001
$ws.Range("H86").Value = 1543.381
$ws.Range("I86").Value = 1444.4375
$ws.Range("J86").Value = 1860
$ws.Range("K86").Value = 1444.4375
$ws.Range("L86").Value = 1860
$ws.Range("M86").Value = -321.4375
$ws.Range("N86").Value = -4106
$ws.Range("H89").Value = 1543.381
$ws.Range("I89").Value = 1444.4375
$ws.Range("J89").Value = 1860
$ws.Range("K89").Value = 7222.1875
$ws.Range("L89").Value = 9300
$ws.Range("M89").Value = -1606.1875
$ws.Range("N89").Value = -20532
$ws.Range("H136").Value = 70400.60000000001
$ws.Range("J136").Value = 70400.60000000001
$ws.Range("L136").Value = 70400.60000000001
$ws.Range("N136").Value = -80600.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1388.1333
$ws.Range("I94").Value = 957
$ws.Range("J94").Value = 1544.909
$ws.Range("K94").Value = 957
$ws.Range("L94").Value = 1544.909
$ws.Range("M94").Value = -506
$ws.Range("N94").Value = -2446.909
$ws.Range("H122").Value = 8936973
$ws.Range("I122").Value = 13900869
$ws.Range("K122").Value = 41702607
$ws.Range("M122").Value = -41700157
$ws.Range("H132").Value = 3129.742
$ws.Range("I132").Value = 2846.4
$ws.Range("J132").Value = 3644.9092
$ws.Range("K132").Value = 8539.200000000001
$ws.Range("L132").Value = 10934.7276
$ws.Range("M132").Value = -6009.200000000001
$ws.Range("N132").Value = -15994.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 18926.273
$ws.Range("I70").Value = 26927
$ws.Range("J70").Value = 4925
$ws.Range("K70").Value = 80781
$ws.Range("L70").Value = 14775
$ws.Range("M70").Value = -80466
$ws.Range("N70").Value = -15405
$ws.Range("H73").Value = 18926.273
$ws.Range("I73").Value = 26927
$ws.Range("J73").Value = 4925
$ws.Range("K73").Value = 80781
$ws.Range("L73").Value = 14775
$ws.Range("M73").Value = -79689
$ws.Range("N73").Value = -16959
$ws.Range("H98").Value = 300
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H131").Value = 5284.174
$ws.Range("I131").Value = 7671.4287
$ws.Range("J131").Value = 4855.6924
$ws.Range("K131").Value = 23014.2861
$ws.Range("L131").Value = 14567.0772
$ws.Range("M131").Value = -17974.2861
$ws.Range("N131").Value = -24647.0772
$ws.Range("H140").Value = 2644.0833
$ws.Range("I140").Value = 2466.125
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 7398.375
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -2218.375
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3319
$ws.Range("I132").Value = 2945.6667
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8837.000100000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6307.000100000001
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2869.6667
$ws.Range("I40").Value = 2808.6667
$ws.Range("J40").Value = 2991.6667
$ws.Range("K40").Value = 2808.6667
$ws.Range("L40").Value = 2991.6667
$ws.Range("M40").Value = -2672.6667
$ws.Range("N40").Value = -3263.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2246.261
$ws.Range("I126").Value = 2348.9092
$ws.Range("J126").Value = 2152.1667
$ws.Range("K126").Value = 7046.7276
$ws.Range("L126").Value = 6456.500100000001
$ws.Range("M126").Value = -4576.7276
$ws.Range("N126").Value = -11396.5001
$ws.Range("H136").Value = 10617.289
$ws.Range("I136").Value = 13398.862
$ws.Range("J136").Value = 1654.4445
$ws.Range("K136").Value = 40196.586
$ws.Range("L136").Value = 4963.333500000001
$ws.Range("M136").Value = -37646.586
$ws.Range("N136").Value = -10063.3335
